$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A55").Value = 2889
$ws.Range("B55").Value = "國票金"
$ws.Range("C55").Value = 0

$ws.Range("A55:C55").Style = $ws.Range("A54:C54").Style

$ws.Range("C55").Select()
